$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "OR" operator symbol cell from "'" to "|"
# (leading apostrophe keeps the existing quote-prefix text style on the cell)
$ws.Range("B6").Value = "'|"

# Update selection to match the recorded view state
$ws.Range("G6").Select()
